# Ajout de tous les données et mise à jour des CNE
# Update the student identifier numbers in column A (rows 2-49): old values
# 1..48 become 330..377 (i.e. +329), reflecting the continuation of the
# roster numbering from a previous batch of students.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 49; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value2 = $cell.Value2 + 329
}

# Update the sheet's active selection / scroll position to match where the
# author ended up after entering the data (near the bottom of the used
# range).
$excel.ActiveWindow.ScrollRow = 48
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("B54").Select()
